$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the date-like text columns stay plain text (as in the source file)
# instead of being auto-converted to date serials.
$ws.Range("Y7:Y8").NumberFormat = "@"
$ws.Range("AA7:AA8").NumberFormat = "@"

# Swap the Id, coordinates, dates and observer between row 7 and row 8
$ws.Range("A7").Value = 107682310
$ws.Range("Q7").Value = 357048.2525418315
$ws.Range("R7").Value = 6425584.133382582
$ws.Range("Y7").Value = "2022-05-22"
$ws.Range("AA7").Value = "2022-05-22"
$ws.Range("AX7").Value = "Olle Kvarnbäck"

$ws.Range("A8").Value = 107683517
$ws.Range("Q8").Value = 357411.1398782768
$ws.Range("R8").Value = 6425417.655266645
$ws.Range("Y8").Value = "2022-03-18"
$ws.Range("AA8").Value = "2022-03-18"
$ws.Range("AX8").Value = "Erik Edvardsson"
